$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells keep a text type, matching the
# original inline-string cells. Without this, numeric-looking values
# like "1.002" would be auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.591.17"
$ws.Range("E2").Value = "  -2.47%  "
$ws.Range("D3").Value = "1.841.55"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "314.54"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.4253"
$ws.Range("E7").Value = "  -3.82%  "
$ws.Range("D8").Value = "0.3642"
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("D9").Value = "45.65"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").Value = "0.07272"
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("D11").Value = "0.8918"
$ws.Range("E11").Value = "  -5.07%  "
$ws.Range("D12").Value = "20.58"
$ws.Range("E12").Value = "  -4.03%  "
$ws.Range("D13").Value = "1.885.63"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "5.373"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").Value = "6.564"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").Value = "0.06870"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "78.37"
$ws.Range("E18").Value = "  -4.64%  "
$ws.Range("D19").Value = "0.000008856"
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "15.53"
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("D22").Value = "27.569.04"
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("D23").Value = "4.983"
$ws.Range("D24").Value = "10.53"
$ws.Range("E24").Value = "  -3.18%  "
$ws.Range("D25").Value = "2.074.43"
$ws.Range("E25").Value = "  -2.77%  "
$ws.Range("D26").Value = "2.040"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").Value = "155.26"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").Value = "18.36"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").Value = "5.206"
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("D30").Value = "116.81"
$ws.Range("E30").Value = "  +2.61%  "
$ws.Range("D31").Value = "1.819"
$ws.Range("E31").Value = "  +5.31%  "
$ws.Range("D32").Value = "0.08893"
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("D33").Value = "0.7771"
$ws.Range("E33").Value = "  -2.93%  "
$ws.Range("D34").Value = "4.557"
$ws.Range("E34").Value = "  -6.37%  "
$ws.Range("D35").Value = "2.976"
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("D36").Value = "1.104"
$ws.Range("E36").Value = "  -6.32%  "
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "0.05406"
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("D39").Value = "1.093"
$ws.Range("E39").Value = "  -3.07%  "
$ws.Range("D40").Value = "0.01917"
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("D41").Value = "2.759"
$ws.Range("E41").Value = "  -8.32%  "
$ws.Range("D42").Value = "6.827"
$ws.Range("E42").Value = "  -4.16%  "
$ws.Range("D43").Value = "0.5058"
$ws.Range("E43").Value = "  -3.92%  "
$ws.Range("D44").Value = "0.1647"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").Value = "8.210"
$ws.Range("E45").Value = "  -6.29%  "
$ws.Range("D46").Value = "0.06619"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("D47").Value = "10.32"
$ws.Range("E47").Value = "  -2.91%  "
$ws.Range("D48").Value = "0.4687"
$ws.Range("E48").Value = "  -4.00%  "
$ws.Range("D49").Value = "105.06"
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("D50").Value = "1.003"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").Value = "1.630"
$ws.Range("E51").Value = "  -2.98%  "
